# Apply cryptos list update (prices + 1h volume %) for Wed Jul 10 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.535.58"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "3.086.47"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "522.67"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.95"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.00%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.439"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.10%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  -0.01%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.385"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("D12").Value = "3.613.15"
$ws.Range("E13").Value = "  +1.00%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.67"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.61%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000166"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "58.544.90"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "3.082.52"
$ws.Range("E17").Value = "  -0.43%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.89"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "8.11"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "342.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +0.00%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.505"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.75"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "0.0₃0916"
$ws.Range("E27").Value = "  -1.84%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.59"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  +1.55%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "21.00"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +2.25%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "154.54"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.05"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.42%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "26.92"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("E37").Value = "  +5.72%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0678"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "3.127.39"
$ws.Range("E39").Value = "  -0.28%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.89"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.75"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.49"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.88%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.667"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "2.271.42"
$ws.Range("E45").Value = "  -0.94%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0257"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.06%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "20.74"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.12%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.956"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.54%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.744"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +7.84%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "265.94"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +9.51%  "
